# The paragraph originally held the "<id>p146v_1</id>" text split across
# three separate runs (the "<id>" / "p146v_1" / "</id>" pieces had been
# typed/pasted in separately, each with its own run formatting). The edit
# collapses that back into a single run so the whole tag+value is one
# contiguous piece of text, carried on the formatting of the first
# ("<id>") run (Courier New, color 7f6000, size 9pt).
#
# Doing a Find/Replace (wdReplaceAll) over the exact visible text re-types
# it as a single run using the formatting in effect at the start of the
# match - i.e. the opening run's formatting - which merges the three runs
# into the one target run and leaves every other paragraph untouched.

$d = $word.ActiveDocument

$targetText = "<id>p146v_1</id>"

$d.Content.Find.Execute($targetText, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $targetText, 2)
